$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.012.43'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '1.558.62'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.44'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.04'
$ws.Range("E8").Value = '  +3.06%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").Value = '1.785.32'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '1.562.85'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("E15").Value = '  +1.86%  '
$ws.Range("D16").Value = '27.034.90'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.92'
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.16'
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.33'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.05'
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.25'
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.28'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.103'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '1.444.43'
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("E34").Value = '  +4.79%  '
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.968'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.523'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.70'
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("E43").Value = '  +3.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.987'
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").Value = '1.697.21'
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.75'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("D50").Value = '0.0₇0987'
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0963'
$ws.Range("E51").Value = '  +2.41%  '
